$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 329.25
$ws.Range("J4").Value = 440
$ws.Range("L4").Value = 440
$ws.Range("N4").Value = -668

$ws.Range("H33").Value = 451.06668
$ws.Range("I33").Value = 85.44444
$ws.Range("J33").Value = 999.5
$ws.Range("K33").Value = 85.44444
$ws.Range("L33").Value = 999.5
$ws.Range("M33").Value = 143.55556
$ws.Range("N33").Value = -1457.5

$ws.Range("H96").Value = 6494517.5
$ws.Range("I96").Value = 14286061
$ws.Range("K96").Value = 42858183
$ws.Range("M96").Value = -42856810

$ws.Range("H100").Value = 1385.2858
$ws.Range("I100").Value = 1290.4546
$ws.Range("J100").Value = 1733
$ws.Range("K100").Value = 1290.4546
$ws.Range("L100").Value = 1733
$ws.Range("M100").Value = -749.4546
$ws.Range("N100").Value = -2815

$ws.Range("H132").Value = 16904.436
$ws.Range("I132").Value = 18205.572
$ws.Range("K132").Value = 54616.716
$ws.Range("M132").Value = -52086.716

$ws.Range("H134").Value = 89999.664
$ws.Range("J134").Value = 89999.664
$ws.Range("L134").Value = 89999.664
$ws.Range("N134").Value = -100139.664

$ws.Range("H135").Value = 1146.0416
$ws.Range("I135").Value = 551.1111
$ws.Range("J135").Value = 2930.8333
$ws.Range("K135").Value = 4959.9999
$ws.Range("L135").Value = 26377.4997
$ws.Range("M135").Value = -2424.9999
$ws.Range("N135").Value = -31447.4997

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 4642.643
$ws.Range("I63").Value = 2500
$ws.Range("K63").Value = 2500
$ws.Range("M63").Value = -1814

$ws.Range("H66").Value = 4642.643
$ws.Range("I66").Value = 2500
$ws.Range("K66").Value = 12500
$ws.Range("M66").Value = -9068

$ws.Range("H74").Value = 255052.95
$ws.Range("I74").Value = 334066.88
$ws.Range("J74").Value = 18011.166
$ws.Range("K74").Value = 334066.88
$ws.Range("L74").Value = 18011.166
$ws.Range("M74").Value = -333192.88
$ws.Range("N74").Value = -19759.166

$ws.Range("H77").Value = 255052.95
$ws.Range("I77").Value = 334066.88
$ws.Range("J77").Value = 18011.166
$ws.Range("K77").Value = 1670334.4
$ws.Range("L77").Value = 90055.83
$ws.Range("M77").Value = -1665966.4
$ws.Range("N77").Value = -98791.83

$ws.Range("H103").Value = 77777
$ws.Range("J103").Value = 77777
$ws.Range("L103").Value = 77777
$ws.Range("N103").Value = -80121

$ws.Range("H132").Value = 1576.561
$ws.Range("I132").Value = 1380.5151
$ws.Range("J132").Value = 2385.25
$ws.Range("K132").Value = 4141.5453
$ws.Range("L132").Value = 7155.75
$ws.Range("M132").Value = -1611.5453
$ws.Range("N132").Value = -12215.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H24").Value = 0
$ws.Range("I24").Value = 0
$ws.Range("K24").Value = 0
$ws.Range("M24").Value = ""

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 635
$ws.Range("J19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("N19").Value = ""

$ws.Range("H24").Value = 635
$ws.Range("J24").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("N24").Value = ""

$ws.Range("H31").Value = 10002900
$ws.Range("I31").Value = 33334332
$ws.Range("K31").Value = 33334332
$ws.Range("M31").Value = -33334037

$ws.Range("H34").Value = 10002900
$ws.Range("I34").Value = 33334332
$ws.Range("K34").Value = 33334332
$ws.Range("M34").Value = -33334130

$ws.Range("H58").Value = 14057.808
$ws.Range("I58").Value = 1366.25
$ws.Range("K58").Value = 1366.25
$ws.Range("M58").Value = -1163.25

$ws.Range("H86").Value = 61849.355
$ws.Range("I86").Value = 75789.3
$ws.Range("J86").Value = 26999.5
$ws.Range("K86").Value = 75789.3
$ws.Range("L86").Value = 26999.5
$ws.Range("M86").Value = -74666.3
$ws.Range("N86").Value = -29245.5

$ws.Range("H89").Value = 61849.355
$ws.Range("I89").Value = 75789.3
$ws.Range("J89").Value = 26999.5
$ws.Range("K89").Value = 378946.5
$ws.Range("L89").Value = 134997.5
$ws.Range("M89").Value = -373330.5
$ws.Range("N89").Value = -146229.5

$ws.Range("H132").Value = 63600.75
$ws.Range("J132").Value = 4904
$ws.Range("L132").Value = 14712
$ws.Range("N132").Value = -19772

$ws.Range("H134").Value = 1682.7188
$ws.Range("I134").Value = 1419.5358
$ws.Range("J134").Value = 3525
$ws.Range("K134").Value = 4258.607400000001
$ws.Range("L134").Value = 10575
$ws.Range("M134").Value = -1723.607400000001
$ws.Range("N134").Value = -15645

$ws.Range("H136").Value = 14057.808
$ws.Range("I136").Value = 1366.25
$ws.Range("K136").Value = 4098.75
$ws.Range("M136").Value = -1548.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 257.33334
$ws.Range("I17").Value = 202.66667
$ws.Range("J17").Value = 366.66666
$ws.Range("K17").Value = 608.00001
$ws.Range("L17").Value = 1099.99998
$ws.Range("M17").Value = -439.00001
$ws.Range("N17").Value = -1437.99998

$ws.Range("H26").Value = 143.9
$ws.Range("I26").Value = 152.85715
$ws.Range("J26").Value = 123
$ws.Range("K26").Value = 458.57145
$ws.Range("L26").Value = 369
$ws.Range("M26").Value = -170.57145
$ws.Range("N26").Value = -945

$ws.Range("H87").Value = 13750
$ws.Range("I87").Value = 20000
$ws.Range("J87").Value = 13478.261
$ws.Range("K87").Value = 60000
$ws.Range("L87").Value = 40434.783
$ws.Range("M87").Value = -58752
$ws.Range("N87").Value = -42930.783

$ws.Range("H90").Value = 13750
$ws.Range("I90").Value = 20000
$ws.Range("J90").Value = 13478.261
$ws.Range("K90").Value = 180000
$ws.Range("L90").Value = 121304.349
$ws.Range("M90").Value = -173760
$ws.Range("N90").Value = -133784.349

$ws.Range("H134").Value = 720.25
$ws.Range("I134").Value = 328.18182
$ws.Range("K134").Value = 984.54546
$ws.Range("M134").Value = 4085.45454

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H29").Value = 17194.666
$ws.Range("I29").Value = 1002
$ws.Range("J29").Value = 20433.2
$ws.Range("K29").Value = 1002
$ws.Range("L29").Value = 20433.2
$ws.Range("M29").Value = -712
$ws.Range("N29").Value = -21013.2

$ws.Range("H97").Value = 827.4
$ws.Range("I97").Value = 860.5
$ws.Range("K97").Value = 860.5
$ws.Range("M97").Value = -364.5

$ws.Range("H134").Value = 41220.668
$ws.Range("J134").Value = 39665
$ws.Range("L134").Value = 118995
$ws.Range("N134").Value = -124065

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 719.0769
$ws.Range("I16").Value = 961.2222
$ws.Range("J16").Value = 174.25
$ws.Range("K16").Value = 961.2222
$ws.Range("L16").Value = 174.25
$ws.Range("M16").Value = -791.2222
$ws.Range("N16").Value = -514.25

$ws.Range("H46").Value = 3469.1365
$ws.Range("I46").Value = 872.4
$ws.Range("K46").Value = 872.4
$ws.Range("M46").Value = -684.4

$ws.Range("H68").Value = 3868.1052
$ws.Range("I68").Value = 3055.4443
$ws.Range("J68").Value = 4599.5
$ws.Range("K68").Value = 3055.4443
$ws.Range("L68").Value = 4599.5
$ws.Range("M68").Value = -2306.4443
$ws.Range("N68").Value = -6097.5

$ws.Range("H71").Value = 3868.1052
$ws.Range("I71").Value = 3055.4443
$ws.Range("J71").Value = 4599.5
$ws.Range("K71").Value = 15277.2215
$ws.Range("L71").Value = 22997.5
$ws.Range("M71").Value = -11533.2215
$ws.Range("N71").Value = -30485.5

$ws.Range("H111").Value = 59998.5
$ws.Range("J111").Value = 59998.5
$ws.Range("L111").Value = 59998.5
$ws.Range("N111").Value = -68178.5

$ws.Range("H122").Value = 3362.2307
$ws.Range("I122").Value = 2718.9678
$ws.Range("K122").Value = 8156.903399999999
$ws.Range("M122").Value = -5706.903399999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 13633
$ws.Range("J5").Value = 13633
$ws.Range("L5").Value = 13633
$ws.Range("N5").Value = -13857

$ws.Range("H104").Value = 21496
$ws.Range("J104").Value = 21496
$ws.Range("L104").Value = 21496
$ws.Range("N104").Value = -28484
